$d = $word.ActiveDocument

# --- 1. Extend the first paragraph with three new runs of text -------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertAfter("This system is based on the assumption that detention will ")
$p1.Range.InsertAfter("be the central focus")
$p1.Range.InsertAfter("; Expulsions, dropouts, etc. will not be factored into this program and counted as a separate process.")

# --- 2. Move the <w:lastRenderedPageBreak/> marker ---------------------------
# It currently sits on the "For what reason..." run; it should instead sit
# on the "What students have demerits" run right before it.
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

$pBreakTarget = $d.Paragraphs.Item(24)
$breakXml = "<w:p w14:paraId='23816348' w14:textId='77777777' w:rsidR='00DA0DC7' w:rsidRDefault='00DA0DC7' w:rsidP='00DA0DC7' $wNs><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='5'/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>What students have demerits</w:t></w:r></w:p>"
$pBreakTarget.Range.InsertXML($breakXml)

$pBreakSource = $d.Paragraphs.Item(25)
$noBreakXml = "<w:p w14:paraId='35CA96F3' w14:textId='77777777' w:rsidR='00DA0DC7' w:rsidRDefault='00DA0DC7' w:rsidP='00DA0DC7' $wNs><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='5'/></w:numPr></w:pPr><w:r><w:t>For what reason the demerit was issued, by category</w:t></w:r></w:p>"
$pBreakSource.Range.InsertXML($noBreakXml)

Write-Host "Done. Paragraph 1 text:" $p1.Range.Text
Write-Host "Paragraph count:" $d.Paragraphs.Count
